# Week 1 "How to write computer algorithms" slide - fix the wording of the
# last bullet point in the "TextBox 8" shape (notes/body textbox), as per
# commit "Edited week 1 object notes".
#
# Before: "Write some computer code that has encodes this set of sequences
#          in a format the CPU can understand."
# After:  "Write some computer code that encodes this set of instructions
#          in a format the CPU can understand."

$p = $ppt.ActivePresentation

$oldText = "Write some computer code that has encodes this set of sequences in a format the CPU can understand."
$newText = "Write some computer code that encodes this set of instructions in a format the CPU can understand."

$targetSlide = $null
$targetShape = $null

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -like "*$oldText*") {
                $targetSlide = $slide
                $targetShape = $shape
                break
            }
        }
    }
    if ($targetShape -ne $null) {
        break
    }
}

$tr = $targetShape.TextFrame.TextRange
$paraCount = $tr.Paragraphs().Count
for ($k = 1; $k -le $paraCount; $k++) {
    $para = $tr.Paragraphs($k)
    if ($para.Text -eq $oldText) {
        $run = $para.Runs(1)
        $run.Text = $newText
        break
    }
}
